$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B8").Value = "Trinity Omics (Acedemic Co-director and Committee Chair)"
$ws.Range("B9").Value = "TTMI Data Science Core"
$ws.Range("B10").Value = "Trinity Single Cell Omics Analysis Workshop"

$ws.Range("C16").Select()
